# Generate Report for Handoff
#
# The "a4a6d581-01fa-4150-9eac-129ac58323aa.md" file has finished translation
# and moved to "Ready for handoff" status. Update its row on the Overview
# sheet as well as the per-locale (zh-cn / de-de) detail sheets with the new
# status, priority, and handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the a4a6d581 file -----------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-08-27 10:15:21"

# --- zh-cn sheet: row 3 is the a4a6d581 file ---------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-27 10:15:17"

# --- de-de sheet: row 3 is the a4a6d581 file ---------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-27 10:15:21"

# --- Re-fit the "Status" columns now that they hold longer text --------
$ovw.Columns("E:F").AutoFit()
$zhcn.Columns("C:C").AutoFit()
$dede.Columns("C:C").AutoFit()
